# Replace the old Nestle product-URL list with the new Aptamil product-URL
# list (16 rows instead of 20) and trim the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://www.uyyaala.com/products/nutricia-aptamil-first-infant-milk-from-birth-0-6-months-800g",
    "https://www.uyyaala.com/products/nutricia-aptamil-pronutra-advanced-milk-formula-to-support-baby-immune-system-stage-2-6-12-months",
    "https://www.uyyaala.com/products/aptamil-lactose-free-first-infant-milk-from-birth-400g",
    "https://www.uyyaala.com/products/nutricia-aptamil-organic-milk-formula-to-support-baby-immune-system-stage-2-6-12-months-800g",
    "https://www.uyyaala.com/products/aptamil-toddler-milk-stage-3-1-2-years-800g",
    "https://www.uyyaala.com/products/nutricia-aptamil-gold-infant-formula-stage-1-from-birth-0-6-months-400g-tin-pack",
    "https://www.uyyaala.com/products/nutricia-aptamil-gold-infant-formula-stage-1-from-birth-0-6-months-400g-refill-pack",
    "https://www.uyyaala.com/products/nutricia-aptamil-gold-infant-formula-stage-2-6-12-months-400g-refill-pack",
    "https://www.uyyaala.com/products/nutricia-aptamil-gold-infant-formula-stage-3-12-months-400g-refill-pack",
    "https://www.uyyaala.com/products/nutricia-aptamil-aptamil-pepti-infant-formula-0-to-12-months-400g",
    "https://www.uyyaala.com/products/nutricia-aptamil-comfort-milk-formula-for-dietary-management-of-colic-constipation-0-12-months",
    "https://www.uyyaala.com/products/nutricia-aptamil-anti-reflux-milk-formula-for-dietary-management-of-reflux-regurgitation-0-12-months-800g",
    "https://www.uyyaala.com/products/nutricia-aptamil-organic-first-infant-milk-from-birth-0-6-months-800g",
    "https://www.uyyaala.com/products/nutricia-aptamil-organic-toddler-milk-stage-3-1-2-years-800g",
    "https://www.uyyaala.com/products/nutricia-aptamil-advanced-3-toddler-milk-substitute-800g-1-3-years",
    "https://www.uyyaala.com/products/nutricia-aptamil-first-infant-milk-from-birth-0-6-months-refill-pack-1-2kg-2-x-600g"
)

# Drop the five rows (18-22) that are no longer needed now that the list has
# shrunk from 20 to 16 entries.
$ws.Range("A18:B22").EntireRow.Delete() | Out-Null

# Overwrite B2:B17 in place with the new URLs (column A keeps its existing
# 0-15 index values and styling).
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $urls[$i]
}

# Header label is unchanged ("urls") but re-assert it for safety.
$ws.Range("B1").Value = "urls"

# Match the saved selection shown in the target workbook.
$ws.Range("C1").Select() | Out-Null
